# Auto-applies the cell-level text updates described by the diff.
# Values are written with NumberFormat forced to Text ("@") so that
# numeric-looking strings (e.g. "1.00", "28.053.51") are preserved
# verbatim as text, matching the workbook's inlineStr-of-text convention,
# then the style is reset to "Normal" so no stray text-format style
# sticks to the cell (keeps cell formatting identical to before).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($ws, [string]$cellRef, [string]$text)
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-CellText $ws 'D2' '28.053.51'
Set-CellText $ws 'E2' '  +3.66%  '
Set-CellText $ws 'D3' '1.726.70'
Set-CellText $ws 'E3' '  +3.03%  '
Set-CellText $ws 'E4' '  -0.11%  '
Set-CellText $ws 'D5' '218.76'
Set-CellText $ws 'E5' '  +1.55%  '
Set-CellText $ws 'E7' '  -0.11%  '
Set-CellText $ws 'E8' '  +13.44%  '
Set-CellText $ws 'E9' '  +3.42%  '
Set-CellText $ws 'E10' '  +1.96%  '
Set-CellText $ws 'D11' '0.0902'
Set-CellText $ws 'E11' '  +2.11%  '
Set-CellText $ws 'D12' '1.970.97'
Set-CellText $ws 'E12' '  +3.07%  '
Set-CellText $ws 'D13' '1.733.36'
Set-CellText $ws 'E13' '  +3.41%  '
Set-CellText $ws 'E14' '  +3.51%  '
Set-CellText $ws 'E15' '  +5.59%  '
Set-CellText $ws 'D16' '67.68'
Set-CellText $ws 'E16' '  +2.58%  '
Set-CellText $ws 'D17' '28.037.62'
Set-CellText $ws 'E17' '  +3.64%  '
Set-CellText $ws 'D18' '243.43'
Set-CellText $ws 'E18' '  +2.40%  '
Set-CellText $ws 'E19' '  +1.92%  '
Set-CellText $ws 'E20' '  -3.31%  '
Set-CellText $ws 'D21' '1.00'
Set-CellText $ws 'E22' '  +3.86%  '
Set-CellText $ws 'E23' '  +4.41%  '
Set-CellText $ws 'D24' '2.14'
Set-CellText $ws 'E24' '  -0.14%  '
Set-CellText $ws 'D25' '149.04'
Set-CellText $ws 'E25' '  +1.63%  '
Set-CellText $ws 'E26' '  +4.52%  '
Set-CellText $ws 'E28' '  +1.94%  '
Set-CellText $ws 'E29' '  -0.12%  '
Set-CellText $ws 'D30' '0.0510'
Set-CellText $ws 'E30' '  +2.55%  '
Set-CellText $ws 'D31' '1.20'
Set-CellText $ws 'E31' '  +2.21%  '
Set-CellText $ws 'B33' 'Maker'
Set-CellText $ws 'C33' 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-CellText $ws 'D33' '1.494.75'
Set-CellText $ws 'E33' '  -3.59%  '
Set-CellText $ws 'B34' 'InternetComputer(DFINITY)'
Set-CellText $ws 'C34' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-CellText $ws 'D34' '3.28'
Set-CellText $ws 'E34' '  +2.76%  '
Set-CellText $ws 'E35' '  -2.41%  '
Set-CellText $ws 'D36' '0.955'
Set-CellText $ws 'E36' '  +3.43%  '
Set-CellText $ws 'D37' '0.608'
Set-CellText $ws 'E37' '  +1.51%  '
Set-CellText $ws 'E38' '  +0.58%  '
Set-CellText $ws 'E39' '  +0.64%  '
Set-CellText $ws 'E40' '  +0.82%  '
Set-CellText $ws 'D41' '70.87'
Set-CellText $ws 'E41' '  +4.81%  '
Set-CellText $ws 'E42' '  +4.20%  '
Set-CellText $ws 'E43' '  -0.09%  '
Set-CellText $ws 'D44' '2.31'
Set-CellText $ws 'E44' '  +2.20%  '
Set-CellText $ws 'D45' '1.874.64'
Set-CellText $ws 'E45' '  +2.86%  '
Set-CellText $ws 'E46' '  +1.75%  '
Set-CellText $ws 'E47' '  +12.16%  '
Set-CellText $ws 'D48' '91.23'
Set-CellText $ws 'E48' '  +0.63%  '
Set-CellText $ws 'E49' '  +3.57%  '
Set-CellText $ws 'E50' '  +0.97%  '
Set-CellText $ws 'D51' '8.16'
Set-CellText $ws 'E51' '  +1.21%  '
